# Auto-generated edit script applying numeric updates to Yojimbo_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 53012.21
$ws.Range("I33").Value = 66986.664
$ws.Range("J33").Value = 608
$ws.Range("K33").Value = 66986.664
$ws.Range("L33").Value = 608
$ws.Range("M33").Value = -66757.664
$ws.Range("N33").Value = -1066

# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 2947.7058
$ws.Range("I113").Value = 2341.1
$ws.Range("J113").Value = 3814.2856
$ws.Range("K113").Value = 2341.1
$ws.Range("L113").Value = 3814.2856
$ws.Range("M113").Value = 912.9000000000001
$ws.Range("N113").Value = -10322.2856

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 863.35
$ws.Range("I137").Value = 806.94116
$ws.Range("J137").Value = 1183
$ws.Range("K137").Value = 2420.82348
$ws.Range("L137").Value = 3549
$ws.Range("M137").Value = 129.17652
$ws.Range("N137").Value = -8649

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1029.8125
$ws.Range("I2").Value = 1039.9166
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 1039.9166
$ws.Range("L2").Value = 999.5
$ws.Range("M2").Value = -926.9166
$ws.Range("N2").Value = -1225.5

# Row 44 (Leve Item ID 3861)
$ws.Range("H44").Value = 39650.75
$ws.Range("I44").Value = 7001
$ws.Range("J44").Value = 44315
$ws.Range("K44").Value = 7001
$ws.Range("L44").Value = 44315
$ws.Range("M44").Value = -6513
$ws.Range("N44").Value = -45291

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1504.3793
$ws.Range("I45").Value = 1520.4286
$ws.Range("J45").Value = 1462.25
$ws.Range("K45").Value = 1520.4286
$ws.Range("L45").Value = 1462.25
$ws.Range("M45").Value = -1143.4286
$ws.Range("N45").Value = -2216.25

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 9576.846
$ws.Range("I61").Value = 10224.917
$ws.Range("K61").Value = 10224.917
$ws.Range("M61").Value = -10012.917

# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 2400
$ws.Range("I88").Value = 2260
$ws.Range("J88").Value = 2750
$ws.Range("K88").Value = 2260
$ws.Range("L88").Value = 2750
$ws.Range("M88").Value = -1854
$ws.Range("N88").Value = -3562

# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 2400
$ws.Range("I91").Value = 2260
$ws.Range("J91").Value = 2750
$ws.Range("K91").Value = 2260
$ws.Range("L91").Value = 2750
$ws.Range("M91").Value = -856
$ws.Range("N91").Value = -5558

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1029.8125
$ws.Range("I116").Value = 1039.9166
$ws.Range("J116").Value = 999.5
$ws.Range("K116").Value = 1039.9166
$ws.Range("L116").Value = 999.5
$ws.Range("M116").Value = 1254.0834
$ws.Range("N116").Value = -5587.5

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 9576.846
$ws.Range("I136").Value = 10224.917
$ws.Range("K136").Value = 30674.751
$ws.Range("M136").Value = -28124.751

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1029.8125
$ws.Range("I3").Value = 1039.9166
$ws.Range("J3").Value = 999.5
$ws.Range("K3").Value = 1039.9166
$ws.Range("L3").Value = 999.5
$ws.Range("M3").Value = -925.9166
$ws.Range("N3").Value = -1227.5

# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 1213.3462
$ws.Range("I20").Value = 1125.55
$ws.Range("J20").Value = 1506
$ws.Range("K20").Value = 1125.55
$ws.Range("L20").Value = 1506
$ws.Range("M20").Value = -878.55
$ws.Range("N20").Value = -2000

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 2237.1667
$ws.Range("I107").Value = 2074.3333
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 2074.3333
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -154.3332999999998
$ws.Range("N107").Value = -6240

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 31812.838
$ws.Range("I31").Value = 40321.574
$ws.Range("J31").Value = 3734
$ws.Range("K31").Value = 40321.574
$ws.Range("L31").Value = 3734
$ws.Range("M31").Value = -40026.574
$ws.Range("N31").Value = -4324

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 31812.838
$ws.Range("I34").Value = 40321.574
$ws.Range("J34").Value = 3734
$ws.Range("K34").Value = 40321.574
$ws.Range("L34").Value = 3734
$ws.Range("M34").Value = -40119.574
$ws.Range("N34").Value = -4138

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 1933.3334
$ws.Range("I58").Value = 1900
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1900
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -1697
$ws.Range("N58").Value = -2406

# Row 94 (Leve Item ID 32934)
$ws.Range("H94").Value = 669853.8
$ws.Range("I94").Value = 672170.7
$ws.Range("J94").Value = 667537
$ws.Range("K94").Value = 672170.7
$ws.Range("L94").Value = 667537
$ws.Range("M94").Value = -671719.7
$ws.Range("N94").Value = -668439

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1336.2632
$ws.Range("I107").Value = 1818.25
$ws.Range("J107").Value = 510
$ws.Range("K107").Value = 1818.25
$ws.Range("L107").Value = 510
$ws.Range("M107").Value = 101.75
$ws.Range("N107").Value = -4350

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1891.2941
$ws.Range("I122").Value = 1831
$ws.Range("J122").Value = 1959.125
$ws.Range("K122").Value = 5493
$ws.Range("L122").Value = 5877.375
$ws.Range("M122").Value = -3043
$ws.Range("N122").Value = -10777.375

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 1933.3334
$ws.Range("I136").Value = 1900
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5700
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3150
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("GSM")
# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 704.13043
$ws.Range("I107").Value = 384.42856
$ws.Range("J107").Value = 1201.4445
$ws.Range("K107").Value = 384.42856
$ws.Range("L107").Value = 1201.4445
$ws.Range("M107").Value = 1535.57144
$ws.Range("N107").Value = -5041.4445

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 1738.3077
$ws.Range("I126").Value = 983.3333
$ws.Range("J126").Value = 2385.4285
$ws.Range("K126").Value = 2949.9999
$ws.Range("L126").Value = 7156.2855
$ws.Range("M126").Value = -479.9998999999998
$ws.Range("N126").Value = -12096.2855

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2205.6843
$ws.Range("I7").Value = 2025
$ws.Range("J7").Value = 2515.4285
$ws.Range("K7").Value = 2025
$ws.Range("L7").Value = 2515.4285
$ws.Range("M7").Value = -1913
$ws.Range("N7").Value = -2739.4285

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2265.0435
$ws.Range("I122").Value = 1941.0667
$ws.Range("J122").Value = 2872.5
$ws.Range("K122").Value = 5823.2001
$ws.Range("L122").Value = 8617.5
$ws.Range("M122").Value = -3373.2001
$ws.Range("N122").Value = -13517.5

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2205.6843
$ws.Range("I126").Value = 2025
$ws.Range("J126").Value = 2515.4285
$ws.Range("K126").Value = 6075
$ws.Range("L126").Value = 7546.2855
$ws.Range("M126").Value = -3605
$ws.Range("N126").Value = -12486.2855

# Row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 33036.223
$ws.Range("J133").Value = 33036.223
$ws.Range("L133").Value = 33036.223
$ws.Range("N133").Value = -38096.223

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 8065.2666
$ws.Range("I136").Value = 8983
$ws.Range("J136").Value = 2100
$ws.Range("K136").Value = 26949
$ws.Range("L136").Value = 6300
$ws.Range("M136").Value = -24399
$ws.Range("N136").Value = -11400
